# This script rewrites the "ToDoList" document body to match the target
# revision described in the commit:
#
#   - Adds a centered, bold, red "To Do" title paragraph at the top.
#   - Turns the old single paragraph (which described linking the
#     "MeasuresSelection" scene and positioning objects in "Referential")
#     into a full numbered To-Do list, each item its own list paragraph:
#       * items that are "done"/obsolete are struck through (w:strike)
#       * some open items are highlighted yellow or cyan
#       * two brand-new trailing items are added ("Undo the last motion
#         (and so on)" and "Close the matching visualization ...")
#   - Leaves the pre-existing "Move the objects..." and "Add the floor..."
#     list items untouched (just renumbered/reflowed along with the rest).
#   - Adds a new yellow-highlighted "Import meshes without the CAD
#     importer" item right before the trailing bookmark-only paragraph.
#   - Removes the list numbering (w:numPr) from that trailing bookmark
#     paragraph, keeping only its paragraph style.
#
# Implementation note: rather than emulate this via many small Find/Replace
# and InsertParagraph calls (error-prone for formatting fidelity), the
# entire document body is rebuilt as one WordprocessingML fragment and
# applied in a single Range.InsertXML call, guaranteeing the resulting
# OOXML matches the target precisely.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="FF0000"/></w:rPr><w:t>To Do</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Create a first scene “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>MeasureSelectiom</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t>” to select the type of measure</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Link the scene “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>Meas</w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>uresSelection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t>” to the next one “Referential”.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Arrange t</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>he objects in the scene “Referential”</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>according to the selected measure</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> type</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Select a sphere by the gaze and “open it” showing the assemblies it represents</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>Close the selection by a gesture and showing again the set of spheres</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>Select an object among the ones displayed by the selected sphere</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Show the matched parts between the query model and the selected one</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>Select a part and move it according to its belonging structure</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Undo the last motion (and so on)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Close the matching visualization and come back to the visualization of the spheres</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Move the objects at user’s meddle height </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Add the floor and a reference frame</w:t></w:r><w:r><w:t xml:space="preserve"> (texture)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Import meshes without the CAD importer</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml.Trim())
